# "Generate Report for Handoff"
# Adds two new source files (b3ec0fd3-... and f6975d88-...) to the
# localization-status workbook, in between the existing "In Translation"
# rows and the ".localization-config" row, on all three sheets
# (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$mdBase  = "https://github.com/OpenLocalizationTest/oltest/blob/5cfc7e0c6f054411fdbe77efc0da23c281f8cf21/e2e/"
$cfgUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5cfc7e0c6f054411fdbe77efc0da23c281f8cf21/.localization-config"
$zhBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b4991f39c2ac089b6018867dbcc7deec321535e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/"
$deBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0192e068a06b8712e7b69198457c77eabe051263/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/"

$newMd1 = "b3ec0fd3-35b1-4067-af9a-34e2b43fac17.md"
$newMd2 = "f6975d88-dda7-44f2-b95b-ae10503d9e02.md"

$newXlfZh1 = "b3ec0fd3-35b1-4067-af9a-34e2b43fac17.77d9b25adf661b6dfc778adca76d6b2e90d9a26c.zh-cn.xlf"
$newXlfZh2 = "f6975d88-dda7-44f2-b95b-ae10503d9e02.1c24bc1aa295aeeff30a45fe4c08e74f979c018d.zh-cn.xlf"
$newXlfDe1 = "b3ec0fd3-35b1-4067-af9a-34e2b43fac17.77d9b25adf661b6dfc778adca76d6b2e90d9a26c.de-de.xlf"
$newXlfDe2 = "f6975d88-dda7-44f2-b95b-ae10503d9e02.1c24bc1aa295aeeff30a45fe4c08e74f979c018d.de-de.xlf"

$zhHandoffDtOld = "2016-02-25 05:09:08"
$deHandoffDtOld = "2016-02-25 05:09:21"
$zhHandoffDtNew = "2016-02-25 05:10:43"
$deHandoffDtNew = "2016-02-25 05:10:55"
$neverDt        = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: "Overview" -- File Name / zh-cn / de-de  (rows 2..6, A:C)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Shift the ".localization-config" summary row down from 4 to 6, and
# insert the two new rows in between the existing rows and it.
$ws1.Range("A6").Value = ".localization-config"
$ws1.Range("B6").Value = "Not to be localized"
$ws1.Range("C6").Value = "Not to be localized"

$ws1.Range("A4").Value = $newMd1
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = $newMd2
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), ($mdBase + "2add2525-e912-4266-b07c-eaf8bcc5659f.md"), "", "", "2add2525-e912-4266-b07c-eaf8bcc5659f.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($mdBase + "512b68af-0e4e-4f2d-a028-22a987387925.md"), "", "", "512b68af-0e4e-4f2d-a028-22a987387925.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), ($mdBase + $newMd1), "", "", $newMd1)
$ws1.Hyperlinks.Add($ws1.Range("A5"), ($mdBase + $newMd2), "", "", $newMd2)
$ws1.Hyperlinks.Add($ws1.Range("A6"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn" detail sheet (rows 2..6, columns A-I)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A6").Value = ".localization-config"
$ws2.Range("B6").Value = "Not to be localized"
$ws2.Range("D6").Value = $neverDt
$ws2.Range("G6").Value = $neverDt
$ws2.Range("H6").Value = "Ignored"

$ws2.Range("A4").Value = $newMd1
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = $newXlfZh1
$ws2.Range("D4").Value = $zhHandoffDtNew
$ws2.Range("G4").Value = $neverDt
$ws2.Range("H4").Value = "Include"

$ws2.Range("A5").Value = $newMd2
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = $newXlfZh2
$ws2.Range("D5").Value = $zhHandoffDtNew
$ws2.Range("G5").Value = $neverDt
$ws2.Range("H5").Value = "Include"

$ws2.Range("D2").Value = $zhHandoffDtOld
$ws2.Range("D3").Value = $zhHandoffDtOld

$ws2.Range("D2:D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), ($mdBase + "2add2525-e912-4266-b07c-eaf8bcc5659f.md"), "", "", "2add2525-e912-4266-b07c-eaf8bcc5659f.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), ($zhBase + "2add2525-e912-4266-b07c-eaf8bcc5659f.255ed9faf4ca82e062f858fa5f6828c948e234d7.zh-cn.xlf"), "", "", "2add2525-e912-4266-b07c-eaf8bcc5659f.255ed9faf4ca82e062f858fa5f6828c948e234d7.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($mdBase + "512b68af-0e4e-4f2d-a028-22a987387925.md"), "", "", "512b68af-0e4e-4f2d-a028-22a987387925.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), ($zhBase + "512b68af-0e4e-4f2d-a028-22a987387925.d158c9bfe3799c8f10567a0ab870f4fb071707de.zh-cn.xlf"), "", "", "512b68af-0e4e-4f2d-a028-22a987387925.d158c9bfe3799c8f10567a0ab870f4fb071707de.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), ($mdBase + $newMd1), "", "", $newMd1)
$ws2.Hyperlinks.Add($ws2.Range("C4"), ($zhBase + $newXlfZh1), "", "", $newXlfZh1)
$ws2.Hyperlinks.Add($ws2.Range("A5"), ($mdBase + $newMd2), "", "", $newMd2)
$ws2.Hyperlinks.Add($ws2.Range("C5"), ($zhBase + $newXlfZh2), "", "", $newXlfZh2)
$ws2.Hyperlinks.Add($ws2.Range("A6"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 3: "de-de" detail sheet (rows 2..6, columns A-I)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A6").Value = ".localization-config"
$ws3.Range("B6").Value = "Not to be localized"
$ws3.Range("D6").Value = $neverDt
$ws3.Range("G6").Value = $neverDt
$ws3.Range("H6").Value = "Ignored"

$ws3.Range("A4").Value = $newMd1
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = $newXlfDe1
$ws3.Range("D4").Value = $deHandoffDtNew
$ws3.Range("G4").Value = $neverDt
$ws3.Range("H4").Value = "Include"

$ws3.Range("A5").Value = $newMd2
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = $newXlfDe2
$ws3.Range("D5").Value = $deHandoffDtNew
$ws3.Range("G5").Value = $neverDt
$ws3.Range("H5").Value = "Include"

$ws3.Range("D2").Value = $deHandoffDtOld
$ws3.Range("D3").Value = $deHandoffDtOld

$ws3.Range("D2:D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), ($mdBase + "2add2525-e912-4266-b07c-eaf8bcc5659f.md"), "", "", "2add2525-e912-4266-b07c-eaf8bcc5659f.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), ($deBase + "2add2525-e912-4266-b07c-eaf8bcc5659f.255ed9faf4ca82e062f858fa5f6828c948e234d7.de-de.xlf"), "", "", "2add2525-e912-4266-b07c-eaf8bcc5659f.255ed9faf4ca82e062f858fa5f6828c948e234d7.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($mdBase + "512b68af-0e4e-4f2d-a028-22a987387925.md"), "", "", "512b68af-0e4e-4f2d-a028-22a987387925.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), ($deBase + "512b68af-0e4e-4f2d-a028-22a987387925.d158c9bfe3799c8f10567a0ab870f4fb071707de.de-de.xlf"), "", "", "512b68af-0e4e-4f2d-a028-22a987387925.d158c9bfe3799c8f10567a0ab870f4fb071707de.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), ($mdBase + $newMd1), "", "", $newMd1)
$ws3.Hyperlinks.Add($ws3.Range("C4"), ($deBase + $newXlfDe1), "", "", $newXlfDe1)
$ws3.Hyperlinks.Add($ws3.Range("A5"), ($mdBase + $newMd2), "", "", $newMd2)
$ws3.Hyperlinks.Add($ws3.Range("C5"), ($deBase + $newXlfDe2), "", "", $newXlfDe2)
$ws3.Hyperlinks.Add($ws3.Range("A6"), $cfgUrl, "", "", ".localization-config")

Write-Output "Done"
